$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for the 380 kV case (Case_2_88)
$updates = @{
    "B2" = 14.82786878980364
    "C2" = 13.26615099274944
    "E2" = 13.38990552783956
    "F2" = 16.86991607391233
    "G2" = 23.48114609952249
    "H2" = 13.13029943223797
    "L2" = 9.825991084247685
    "M2" = 14.19150112599157
    "O2" = 19.25379860382921
    "B3" = 14.19808162115057
    "C3" = 13.13080476373325
    "E3" = 13.45949124857107
    "F3" = 15.89584955866808
    "G3" = 23.70662309395977
    "H3" = 13.20307032929184
    "L3" = 9.832305485148286
    "M3" = 14.04287027211824
    "O3" = 19.39196281608329
    "B4" = 13.79670086905288
    "C4" = 13.04804393776718
    "E4" = 13.50447946049354
    "F4" = 15.26997757108491
    "G4" = 23.85614092975618
    "H4" = 13.25038677684129
    "L4" = 9.837493439078253
    "M4" = 13.95197520415265
    "O4" = 19.48224176457315
    "B5" = 13.62964313649275
    "C5" = 13.01443353439404
    "E5" = 13.52338304414193
    "F5" = 15.00819731993403
    "G5" = 23.91983727918764
    "H5" = 13.2703317749575
    "L5" = 9.839937797280397
    "M5" = 13.91505825707945
    "O5" = 19.52039907098326
    "B6" = 13.60169886886438
    "C6" = 13.00886036760821
    "E6" = 13.52655647645526
    "F6" = 14.96433081551593
    "G6" = 23.93058058607107
    "H6" = 13.2736836985908
    "L6" = 9.840363641848706
    "M6" = 13.90893661149962
    "O6" = 19.52681764078637
    "B7" = 13.79446172147215
    "C7" = 13.04759015075574
    "E7" = 13.50473208861759
    "F7" = 15.26647399323137
    "G7" = 23.85698878265794
    "H7" = 13.2506530761154
    "L7" = 9.837525066746082
    "M7" = 13.9514767873073
    "O7" = 19.48275083059031
    "B8" = 14.61387962987321
    "C8" = 13.21943013817943
    "E8" = 13.41343013917414
    "F8" = 16.53996406344768
    "G8" = 23.55658119574514
    "H8" = 13.15484434702624
    "L8" = 9.827896588424807
    "M8" = 14.14019659285239
    "O8" = 19.30030651770968
    "B9" = 16.09665780592045
    "C9" = 13.5578382767969
    "E9" = 13.25226209564852
    "F9" = 19.00274580682531
    "G9" = 23.0561728557312
    "H9" = 12.98784179012228
    "L9" = 9.819388678670943
    "M9" = 14.51177117373209
    "O9" = 18.98580770434948
    "B10" = 17.10193981432359
    "C10" = 13.80562397937308
    "E10" = 13.14464278877724
    "F10" = 20.67494806633232
    "G10" = 22.74374000880306
    "H10" = 12.8778336077112
    "L10" = 9.819421879862462
    "M10" = 14.78381007322585
    "O10" = 18.78121983766327
    "B11" = 17.53962686361628
    "C11" = 13.91780580668235
    "E11" = 13.09800454091681
    "F11" = 21.3917225636224
    "G11" = 22.6138510832658
    "H11" = 12.8305355126026
    "L11" = 9.820791537642689
    "M11" = 14.90696780577779
    "O11" = 18.69391884659272
    "B12" = 17.70245413787701
    "C12" = 13.96017932102222
    "E12" = 13.08067552467967
    "F12" = 21.65686569030329
    "G12" = 22.56644718889351
    "H12" = 12.81301925833625
    "L12" = 9.821503984838586
    "M12" = 14.95348669709925
    "O12" = 18.66169181606399
    "B13" = 17.66751736825356
    "C13" = 13.95105867137417
    "E13" = 13.0843929017708
    "F13" = 21.60004134736742
    "G13" = 22.57657684078264
    "H13" = 12.81677415474651
    "L13" = 9.821341944533877
    "M13" = 14.94347379731577
    "O13" = 18.66859543808839
    "B14" = 17.55308167373394
    "C14" = 13.9212942218971
    "E14" = 13.0965722295466
    "F14" = 21.4136618050453
    "G14" = 22.6099152885127
    "H14" = 12.82908653456368
    "L14" = 9.820846273296381
    "M14" = 14.91079749897808
    "O14" = 18.69125081990036
    "B15" = 17.48260431483276
    "C15" = 13.9030477803203
    "E15" = 13.10407559539897
    "F15" = 21.29868154950795
    "G15" = 22.63056881855407
    "H15" = 12.83667959266565
    "L15" = 9.820567866304376
    "M15" = 14.89076595730574
    "O15" = 18.70523632272646
    "B16" = 17.07293274697975
    "C16" = 13.79827924063917
    "E16" = 13.14773723029809
    "F16" = 20.62722412089977
    "G16" = 22.75247672535561
    "H16" = 12.88097986759549
    "L16" = 9.819359543289677
    "M16" = 14.77574664907968
    "O16" = 18.78704143513895
    "B17" = 16.81651550816649
    "C17" = 13.73384798345178
    "E17" = 13.17511493737661
    "F17" = 20.20408069597325
    "G17" = 22.83041360144163
    "H17" = 12.90885957046198
    "L17" = 9.818964526758982
    "M17" = 14.70501029934856
    "O17" = 18.8387051536924
    "B18" = 16.66718887315354
    "C18" = 13.6967400360095
    "E18" = 13.19108016965359
    "F18" = 19.95656407809801
    "G18" = 22.87639139974359
    "H18" = 12.92515363960021
    "L18" = 9.818864885345066
    "M18" = 14.66427062780964
    "O18" = 18.86896341243712
    "B19" = 16.61631614639342
    "C19" = 13.68416846830051
    "E19" = 13.19652326185708
    "F19" = 19.87204792380568
    "G19" = 22.89215562836772
    "H19" = 12.93071492026885
    "L19" = 9.818853086354679
    "M19" = 14.65046861227314
    "O19" = 18.87930146151693
    "B20" = 16.84400295147319
    "C20" = 13.74071207032453
    "E20" = 13.17217795143345
    "F20" = 20.24955283636154
    "G20" = 22.8219978476047
    "H20" = 12.90586498681312
    "L20" = 9.818993381754803
    "M20" = 14.71254613259287
    "O20" = 18.83314928111539
    "B21" = 17.58677400662108
    "C21" = 13.93003991114498
    "E21" = 13.09298587057482
    "F21" = 21.46857628470577
    "G21" = 22.60007442917082
    "H21" = 12.82545938508671
    "L21" = 9.820986612944433
    "M21" = 14.92039878807429
    "O21" = 18.68457378397797
    "B22" = 18.05519168026181
    "C22" = 14.05313762693899
    "E22" = 13.04316301101555
    "F22" = 22.22866616901552
    "G22" = 22.46543344097961
    "H22" = 12.775209084086
    "L22" = 9.823418328895473
    "M22" = 15.05553841046098
    "O22" = 18.59232174739072
    "B23" = 17.80677351498759
    "C23" = 13.98750629618382
    "E23" = 13.06957796167602
    "F23" = 21.82633154458858
    "G23" = 22.5363348952163
    "H23" = 12.8018182903619
    "L23" = 9.822017519957861
    "M23" = 14.98348694760096
    "O23" = 18.6411136437488
    "B24" = 16.83158181703953
    "C24" = 13.73760901856145
    "E24" = 13.1735050604878
    "F24" = 20.22900810905287
    "G24" = 22.82579896448055
    "H24" = 12.90721801054895
    "L24" = 9.818979939317316
    "M24" = 14.70913940538666
    "O24" = 18.83565935908161
    "B25" = 15.70976964044982
    "C25" = 13.4663145466466
    "E25" = 13.29395971690963
    "F25" = 18.34778573295695
    "G25" = 23.18193461126831
    "H25" = 13.03078945759021
    "L25" = 9.820583571902661
    "M25" = 14.41128392395163
    "O25" = 19.06624603512856
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}

